$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (Version, Date, Contact) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Update the concept value on the original "Include from FSIII" sheet ---
$src = $wb.Worksheets.Item("Include from FSIII")
$src.Range("C2").Value = "ad78224f-b339-462c-9f2c-90b3120605cb"

# --- Duplicate the sheet to create "Include from FSIII 2" with 4 concepts (FBOE) ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Include from FSIII 2"
$new.Range("C2").Value = "E"
